# Refresh the crypto market snapshot (prices + 1h volume deltas) to the latest pull.
# Row 39/38 also reflects an upstream rank swap: Kaspa and Stacks traded places.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '60.874.74'
$ws.Range('E2').Value = '  +0.01%  '

# Row 3
$ws.Range('D3').Value = '2.913.12'
$ws.Range('E3').Value = '  -0.23%  '

# Row 4
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = '  -0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '591.05'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.85%  '

# Row 6
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.81'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.01%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.506'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +0.55%  '

# Row 9
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '6.90'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.98%  '

# Row 10
$ws.Range('E10').Value = '  -0.94%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.440'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.74%  '

# Row 12
$ws.Range('E12').Value = '  -0.25%  '

# Row 13
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '33.54'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  -0.53%  '

# Row 14
$ws.Range('E14').Value = '  -0.12%  '

# Row 15
$ws.Range('D15').Value = '3.394.97'
$ws.Range('E15').Value = '  -0.23%  '

# Row 16
$ws.Range('D16').Value = '60.815.64'
$ws.Range('E16').Value = '  -0.01%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '6.69'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -1.15%  '

# Row 18
$ws.Range('D18').Value = '2.914.06'
$ws.Range('E18').Value = '  -0.11%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '431.19'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.33%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.36'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.38%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.678'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -0.53%  '

# Row 22
$ws.Range('E22').Value = '  -1.77%  '

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '81.63'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +1.14%  '

# Row 24
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.05'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +0.21%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.19'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  -1.13%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.82'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.47%  '

# Row 27
$ws.Range('E27').Value = '  +0.01%  '

# Row 28
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.27'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +4.46%  '

# Row 29
$ws.Range('E29').Value = '  -0.68%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '6.99'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -3.05%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '26.51'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.54%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.108'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.94%  '

# Row 33
$ws.Range('E33').Value = '  +0.02%  '

# Row 34
$ws.Range('D34').Value = '0.0₃0864'
$ws.Range('E34').Value = '  -2.37%  '

# Row 35
$ws.Range('E35').Value = '  -0.24%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.62'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -0.70%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.03'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -0.38%  '

# Row 38
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.98'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -1.38%  '

# Row 39
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.122'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.01%  '

# Row 40
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '8.54'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -1.83%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.286'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -4.22%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '40.72'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  -2.33%  '

# Row 43
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '377.69'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.58%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0343'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -1.72%  '

# Row 45
$ws.Range('D45').Value = '2.693.20'
$ws.Range('E45').Value = '  +0.56%  '

# Row 46
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '133.30'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +1.03%  '

# Row 47
$ws.Range('E47').Value = '  -0.03%  '

# Row 48
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '23.88'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -3.28%  '

# Row 49
$ws.Range('E49').Value = '  -0.69%  '

# Row 50
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '2.00'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -3.44%  '

# Row 51
$ws.Range('E51').Value = '  +0.04%  '
